# Update the "ResumoInscricoes_Superior" data sheet with refreshed counts.
# The source data changed (more registrations/payments came in), so a set of
# individual cells in columns E (Inscritos), F (Pagos) and H (Inscrições
# homologadas) need to be bumped to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    "E2"  = 54
    "F6"  = 3
    "H6"  = 4
    "F15" = 93
    "H15" = 134
    "E16" = 13
    "F16" = 5
    "H16" = 9
    "F17" = 67
    "H17" = 99
    "E18" = 122
    "F18" = 56
    "H18" = 92
    "E19" = 66
    "F19" = 39
    "H19" = 52
    "E23" = 7
    "F23" = 2
    "H23" = 4
    "E24" = 25
    "F24" = 16
    "H24" = 20
    "E25" = 24
    "E26" = 30
    "E27" = 11
    "E35" = 10
    "F35" = 5
    "H35" = 6
    "F36" = 53
    "H36" = 85
    "F37" = 34
    "H37" = 46
    "F40" = 16
    "H40" = 19
    "E41" = 45
    "F41" = 22
    "H41" = 33
    "E42" = 39
    "F42" = 22
    "H42" = 31
    "E43" = 26
    "E44" = 30
    "F44" = 16
    "H44" = 26
    "F46" = 12
    "H46" = 21
    "E48" = 34
    "F48" = 23
    "H48" = 28
    "E49" = 75
    "F49" = 41
    "H49" = 58
    "E50" = 26
    "F50" = 10
    "H50" = 18
    "F52" = 5
    "H52" = 5
    "F60" = 11
    "H60" = 16
    "F61" = 12
    "H61" = 22
    "F62" = 14
    "H62" = 28
    "F63" = 13
    "H63" = 21
    "F64" = 19
    "H64" = 24
    "E65" = 36
    "F65" = 12
    "H65" = 25
    "E66" = 36
    "F66" = 24
    "H66" = 32
    "F69" = 9
    "H69" = 12
    "F70" = 25
    "H70" = 37
    "F72" = 27
    "H72" = 38
    "F73" = 11
    "H73" = 23
    "F74" = 9
    "H74" = 13
    "E75" = 16
    "F75" = 8
    "H75" = 13
    "F76" = 20
    "H76" = 37
    "E79" = 41
    "F79" = 20
    "H79" = 32
    "F82" = 6
    "H82" = 12
    "F83" = 3
    "H83" = 10
    "F87" = 7
    "H87" = 14
    "F88" = 16
    "H88" = 24
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
